# Update Mappings 22 Ontologies
# Adds a new "MS_DEF" column (F) to the mapping sheet, with an empty-list
# placeholder value "[]" for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (21 in this workbook, but computed
# dynamically so the script is robust).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# --- New header cell F1 -------------------------------------------------
# Give F1 the same look (bold / border / centered) as the other header
# cells by copying the formatting from E1, then set its own text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "MS_DEF"

# --- New data column F2:F<lastRow> --------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
